# Natmi following Dr Hou advice
# Update the Sfrp1-Fzd6 LR-pair sheet: add "sCs" sending-cluster rows and
# refresh the computed statistics for every existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Sfrp1", "Fzd6", "ECs",  3, 1, 3.015829666666667, 9.047489000000001, 0.02149774561142444, 0.02149774561142444, 3, 1, 20.545366, 61.636098, 0.8965950288338865, 0.8965950288338865, 61.96132429532467, 557.6519186579221, 0.01927477184633865, 0.01927477184633865),
    @("ECs",  "Sfrp1", "Fzd6", "FAPs", 3, 1, 3.015829666666667, 9.047489000000001, 0.02149774561142444, 0.02149774561142444, 3, 1, 2.154739666666666, 6.464219, 0.09403234157836461, 0.09403234157836463, 6.498327810676778, 58.484950296091, 0.002021483358498252, 0.002021483358498252),
    @("ECs",  "Sfrp1", "Fzd6", "sCs",  3, 1, 3.015829666666667, 9.047489000000001, 0.02149774561142444, 0.02149774561142444, 3, 1, 0.2147726666666667, 0.6443180000000001, 0.009372629587748921, 0.009372629587748921, 0.6477177797224446, 5.829460017502001, 0.0002014904065875362, 0.0002014904065875362),
    @("FAPs", "Sfrp1", "Fzd6", "ECs",  3, 1, 136.9593173333334, 410.8779520000001, 0.976287419574543, 0.976287419574543, 3, 1, 20.545366, 61.636098, 0.8965950288338865, 0.8965950288338865, 2813.879301723478, 25324.9137155113, 0.8753344471035981, 0.8753344471035981),
    @("FAPs", "Sfrp1", "Fzd6", "FAPs", 3, 1, 136.9593173333334, 410.8779520000001, 0.976287419574543, 0.976287419574543, 3, 1, 2.154739666666666, 6.464219, 0.09403234157836461, 0.09403234157836463, 295.1116737777209, 2656.005063999488, 0.0918025921160936, 0.09180259211609361),
    @("FAPs", "Sfrp1", "Fzd6", "sCs",  3, 1, 136.9593173333334, 410.8779520000001, 0.976287419574543, 0.976287419574543, 3, 1, 0.2147726666666667, 0.6443180000000001, 0.009372629587748921, 0.009372629587748921, 29.41511780852623, 264.7360602767361, 0.009150380354851407, 0.009150380354851407),
    @("sCs",  "Sfrp1", "Fzd6", "ECs",  2, 0.6666666666666666, 0.31071, 0.93213, 0.002214834814032607, 0.002214834814032608, 3, 1, 20.545366, 61.636098, 0.8965950288338865, 0.8965950288338865, 6.38365066986, 57.45285602874, 0.001985809883949861, 0.001985809883949861),
    @("sCs",  "Sfrp1", "Fzd6", "FAPs", 2, 0.6666666666666666, 0.31071, 0.93213, 0.002214834814032607, 0.002214834814032608, 3, 1, 2.154739666666666, 6.464219, 0.09403234157836461, 0.09403234157836463, 0.6694991618299999, 6.02549245647, 0.0002082661037727678, 0.0002082661037727678),
    @("sCs",  "Sfrp1", "Fzd6", "sCs",  2, 0.6666666666666666, 0.31071, 0.93213, 0.002214834814032607, 0.002214834814032608, 3, 1, 0.2147726666666667, 0.6443180000000001, 0.009372629587748921, 0.009372629587748921, 0.06673201526, 0.6005881373400001, 0.00002075882630997839, 0.0000207588263099784)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
